$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("articels")

# Row 4 now gets a value (previously empty)
$ws.Range("A4").Value = "/add_vendor"

# New rows 9-13
$ws.Range("A9").Value  = "/add_vendor"
$ws.Range("A10").Value = "/addven"
$ws.Range("A11").Value = "/addven"
$ws.Range("A12").Value = "/addven"
$ws.Range("A13").Value = "/addven"
